$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a text value to a cell without letting Excel
# auto-convert numeric-looking strings (e.g. "320.90") into a
# number. Cells in column D store plain-text price strings, so
# a write like "320.90" must stay the literal text "320.90" and
# not become the number 320.9. Forcing NumberFormat = "@" (Text)
# before the assignment keeps it a string; resetting the style
# back to "Normal" afterwards avoids leaving a stray number-format
# on the cell (it had none before the edit).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "47.200.59"
$ws.Cells.Item(2, 5).Value = "  -0.33%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.485.52"
$ws.Cells.Item(3, 5).Value = "  -0.94%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
Set-TextValue $ws.Cells.Item(5, 4) "320.90"
$ws.Cells.Item(5, 5).Value = "  -1.18%  "
Set-TextValue $ws.Cells.Item(6, 4) "107.55"
$ws.Cells.Item(6, 5).Value = "  +1.21%  "
$ws.Cells.Item(7, 5).Value = "  -0.61%  "
$ws.Cells.Item(8, 5).Value = "  -0.02%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.532"
$ws.Cells.Item(9, 5).Value = "  -1.65%  "
Set-TextValue $ws.Cells.Item(10, 4) "38.43"
$ws.Cells.Item(10, 5).Value = "  +4.54%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.0807"
$ws.Cells.Item(11, 5).Value = "  -1.40%  "
$ws.Cells.Item(12, 5).Value = "  -0.04%  "
Set-TextValue $ws.Cells.Item(13, 4) "18.21"
$ws.Cells.Item(13, 5).Value = "  -1.07%  "
Set-TextValue $ws.Cells.Item(14, 4) "7.09"
$ws.Cells.Item(14, 5).Value = "  -1.26%  "
Set-TextValue $ws.Cells.Item(15, 4) "2.873.24"
$ws.Cells.Item(15, 5).Value = "  -1.05%  "
Set-TextValue $ws.Cells.Item(16, 4) "2.487.91"
$ws.Cells.Item(16, 5).Value = "  -2.66%  "
$ws.Cells.Item(17, 5).Value = "  -0.33%  "
Set-TextValue $ws.Cells.Item(18, 4) "47.113.95"
$ws.Cells.Item(18, 5).Value = "  -0.34%  "
Set-TextValue $ws.Cells.Item(19, 4) "12.71"
$ws.Cells.Item(19, 5).Value = "  -0.75%  "
$ws.Cells.Item(20, 5).Value = "  +1.40%  "
Set-TextValue $ws.Cells.Item(21, 4) "0.0₃0930"
$ws.Cells.Item(21, 5).Value = "  -1.41%  "
Set-TextValue $ws.Cells.Item(22, 4) "2.71"
$ws.Cells.Item(22, 5).Value = "  +13.45%  "
Set-TextValue $ws.Cells.Item(23, 4) "70.19"
$ws.Cells.Item(23, 5).Value = "  -1.12%  "
Set-TextValue $ws.Cells.Item(24, 4) "245.25"
$ws.Cells.Item(24, 5).Value = "  -2.99%  "
$ws.Cells.Item(25, 5).Value = "  +0.06%  "
$ws.Cells.Item(26, 5).Value = "  -0.03%  "
Set-TextValue $ws.Cells.Item(27, 4) "25.62"
$ws.Cells.Item(27, 5).Value = "  -3.09%  "
$ws.Cells.Item(28, 5).Value = "  +0.17%  "
Set-TextValue $ws.Cells.Item(29, 4) "9.99"
$ws.Cells.Item(29, 5).Value = "  +0.14%  "
Set-TextValue $ws.Cells.Item(30, 4) "34.32"
$ws.Cells.Item(30, 5).Value = "  -2.63%  "
$ws.Cells.Item(31, 5).Value = "  -1.61%  "
Set-TextValue $ws.Cells.Item(32, 4) "49.47"
$ws.Cells.Item(32, 5).Value = "  -0.70%  "
Set-TextValue $ws.Cells.Item(33, 4) "20.23"
$ws.Cells.Item(33, 5).Value = "  +1.96%  "
Set-TextValue $ws.Cells.Item(34, 4) "5.32"
$ws.Cells.Item(34, 5).Value = "  -0.19%  "
$ws.Cells.Item(35, 5).Value = "  -0.27%  "
$ws.Cells.Item(36, 5).Value = "  +0.14%  "
Set-TextValue $ws.Cells.Item(37, 4) "1.96"
$ws.Cells.Item(37, 5).Value = "  +0.23%  "
$ws.Cells.Item(38, 5).Value = "  -1.14%  "
$ws.Cells.Item(39, 5).Value = "  -1.48%  "
$ws.Cells.Item(40, 5).Value = "  -0.86%  "
Set-TextValue $ws.Cells.Item(41, 4) "22.26"
$ws.Cells.Item(41, 5).Value = "  +3.00%  "
$ws.Cells.Item(42, 5).Value = "  -0.62%  "
Set-TextValue $ws.Cells.Item(43, 4) "118.79"
$ws.Cells.Item(43, 5).Value = "  -4.33%  "
$ws.Cells.Item(44, 5).Value = "  -1.32%  "
Set-TextValue $ws.Cells.Item(45, 4) "1.982.37"
$ws.Cells.Item(45, 5).Value = "  +0.09%  "
$ws.Cells.Item(46, 5).Value = "  -1.22%  "
Set-TextValue $ws.Cells.Item(47, 4) "1.99"
$ws.Cells.Item(47, 5).Value = "  -6.47%  "
Set-TextValue $ws.Cells.Item(48, 4) "9.02"
$ws.Cells.Item(48, 5).Value = "  -0.72%  "
Set-TextValue $ws.Cells.Item(49, 4) "1.76"
$ws.Cells.Item(49, 5).Value = "  -2.93%  "
Set-TextValue $ws.Cells.Item(50, 4) "5.10"
$ws.Cells.Item(50, 5).Value = "  -6.12%  "
Set-TextValue $ws.Cells.Item(51, 4) "56.50"
$ws.Cells.Item(51, 5).Value = "  +3.01%  "
